# Auto-generated script applying scheduled market-data refresh to Titan_Profits sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4171056.5
$ws.Range("I64").Value = 8931983
$ws.Range("J64").Value = 5245.9375
$ws.Range("K64").Value = 8931983
$ws.Range("L64").Value = 5245.9375
$ws.Range("M64").Value = -8931735
$ws.Range("N64").Value = -5741.9375
$ws.Range("H67").Value = 4171056.5
$ws.Range("I67").Value = 8931983
$ws.Range("J67").Value = 5245.9375
$ws.Range("K67").Value = 8931983
$ws.Range("L67").Value = 5245.9375
$ws.Range("M67").Value = -8931125
$ws.Range("N67").Value = -6961.9375
$ws.Range("H76").Value = 3271079
$ws.Range("I76").Value = 4632504.5
$ws.Range("J76").Value = 3658
$ws.Range("K76").Value = 4632504.5
$ws.Range("L76").Value = 3658
$ws.Range("M76").Value = -4632189.5
$ws.Range("N76").Value = -4288
$ws.Range("H79").Value = 3271079
$ws.Range("I79").Value = 4632504.5
$ws.Range("J79").Value = 3658
$ws.Range("K79").Value = 4632504.5
$ws.Range("L79").Value = 3658
$ws.Range("M79").Value = -4631412.5
$ws.Range("N79").Value = -5842
$ws.Range("H107").Value = 505609.72
$ws.Range("I107").Value = 585374.4399999999
$ws.Range("J107").Value = 433.33334
$ws.Range("K107").Value = 585374.4399999999
$ws.Range("L107").Value = 433.33334
$ws.Range("M107").Value = -583454.4399999999
$ws.Range("N107").Value = -4273.33334
$ws.Range("H112").Value = 14355412
$ws.Range("J112").Value = 14355412
$ws.Range("L112").Value = 43066236
$ws.Range("N112").Value = -43068452
$ws.Range("H137").Value = 47621690
$ws.Range("I137").Value = 71429900
$ws.Range("J137").Value = 5257.5713
$ws.Range("K137").Value = 214289700
$ws.Range("L137").Value = 15772.7139
$ws.Range("M137").Value = -214287150
$ws.Range("N137").Value = -20872.7139

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 200700.2
$ws.Range("I5").Value = 250850.25
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 250850.25
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -250738.25
$ws.Range("N5").Value = -324
$ws.Range("H63").Value = 16900
$ws.Range("I63").Value = 24000
$ws.Range("K63").Value = 24000
$ws.Range("M63").Value = -23314
$ws.Range("H66").Value = 16900
$ws.Range("I66").Value = 24000
$ws.Range("K66").Value = 120000
$ws.Range("M66").Value = -116568
$ws.Range("H74").Value = 6880.391
$ws.Range("I74").Value = 1290.5
$ws.Range("J74").Value = 12978.454
$ws.Range("K74").Value = 1290.5
$ws.Range("L74").Value = 12978.454
$ws.Range("M74").Value = -416.5
$ws.Range("N74").Value = -14726.454
$ws.Range("H77").Value = 6880.391
$ws.Range("I77").Value = 1290.5
$ws.Range("J77").Value = 12978.454
$ws.Range("K77").Value = 6452.5
$ws.Range("L77").Value = 64892.27
$ws.Range("M77").Value = -2084.5
$ws.Range("N77").Value = -73628.26999999999
$ws.Range("H139").Value = 51905
$ws.Range("J139").Value = 51905
$ws.Range("L139").Value = 51905
$ws.Range("N139").Value = -62185

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200700.2
$ws.Range("I4").Value = 250850.25
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 250850.25
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -250735.25
$ws.Range("N4").Value = -330
$ws.Range("H82").Value = 7726.3335
$ws.Range("I82").Value = 7726.3335
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 7726.3335
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -7343.3335
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 7726.3335
$ws.Range("I85").Value = 7726.3335
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 7726.3335
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -6400.3335
$ws.Range("N85").ClearContents()
$ws.Range("H105").Value = 3693.606
$ws.Range("I105").Value = 3608.5908
$ws.Range("J105").Value = 3863.6365
$ws.Range("K105").Value = 3608.5908
$ws.Range("L105").Value = 3863.6365
$ws.Range("M105").Value = -1861.5908
$ws.Range("N105").Value = -7357.636500000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4323.5537
$ws.Range("I31").Value = 2105.2307
$ws.Range("J31").Value = 5802.436
$ws.Range("K31").Value = 2105.2307
$ws.Range("L31").Value = 5802.436
$ws.Range("M31").Value = -1810.2307
$ws.Range("N31").Value = -6392.436
$ws.Range("H34").Value = 4323.5537
$ws.Range("I34").Value = 2105.2307
$ws.Range("J34").Value = 5802.436
$ws.Range("K34").Value = 2105.2307
$ws.Range("L34").Value = 5802.436
$ws.Range("M34").Value = -1903.2307
$ws.Range("N34").Value = -6206.436
$ws.Range("H86").Value = 71430840
$ws.Range("I86").Value = 166667950
$ws.Range("J86").Value = 2999.75
$ws.Range("K86").Value = 166667950
$ws.Range("L86").Value = 2999.75
$ws.Range("M86").Value = -166666827
$ws.Range("N86").Value = -5245.75
$ws.Range("H89").Value = 71430840
$ws.Range("I89").Value = 166667950
$ws.Range("J89").Value = 2999.75
$ws.Range("K89").Value = 833339750
$ws.Range("L89").Value = 14998.75
$ws.Range("M89").Value = -833334134
$ws.Range("N89").Value = -26230.75
$ws.Range("H123").Value = 33000
$ws.Range("J123").Value = 33000
$ws.Range("L123").Value = 33000
$ws.Range("N123").Value = -42800
$ws.Range("H141").Value = 224541.44
$ws.Range("J141").Value = 228618.45
$ws.Range("L141").Value = 228618.45
$ws.Range("N141").Value = -238978.45

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1380.2667
$ws.Range("I132").Value = 1196.4
$ws.Range("K132").Value = 10767.6
$ws.Range("M132").Value = -8237.6

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4829.722
$ws.Range("I70").Value = 4864.4316
$ws.Range("K70").Value = 4864.4316
$ws.Range("M70").Value = -4594.4316
$ws.Range("H73").Value = 4829.722
$ws.Range("I73").Value = 4864.4316
$ws.Range("K73").Value = 4864.4316
$ws.Range("M73").Value = -3928.4316
$ws.Range("H80").Value = 38464436
$ws.Range("I80").Value = 2692.5334
$ws.Range("J80").Value = 90912264
$ws.Range("K80").Value = 2692.5334
$ws.Range("L80").Value = 90912264
$ws.Range("M80").Value = -1694.5334
$ws.Range("N80").Value = -90914260
$ws.Range("H83").Value = 38464436
$ws.Range("I83").Value = 2692.5334
$ws.Range("J83").Value = 90912264
$ws.Range("K83").Value = 13462.667
$ws.Range("L83").Value = 454561320
$ws.Range("M83").Value = -8470.666999999999
$ws.Range("N83").Value = -454571304
$ws.Range("H113").Value = 1800.3334
$ws.Range("I113").Value = 1800.5
$ws.Range("K113").Value = 1800.5
$ws.Range("M113").Value = 369.5
$ws.Range("H122").Value = 2756.2
$ws.Range("I122").Value = 3033.75
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9101.25
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6651.25
$ws.Range("N122").Value = -12400

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 573.1
$ws.Range("I46").Value = 513.5
$ws.Range("J46").Value = 662.5
$ws.Range("K46").Value = 513.5
$ws.Range("L46").Value = 662.5
$ws.Range("M46").Value = -325.5
$ws.Range("N46").Value = -1038.5
$ws.Range("H68").Value = 1920
$ws.Range("I68").Value = 1908.5714
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1908.5714
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1159.5714
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1920
$ws.Range("I71").Value = 1908.5714
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9542.857
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -5798.857
$ws.Range("N71").Value = -17488
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 2479.1428
$ws.Range("I93").Value = 2088.5
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2088.5
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -840.5
$ws.Range("N93").Value = -5496
$ws.Range("H122").Value = 3347.5454
$ws.Range("I122").Value = 2435.6428
$ws.Range("J122").Value = 4019.4736
$ws.Range("K122").Value = 7306.928400000001
$ws.Range("L122").Value = 12058.4208
$ws.Range("M122").Value = -4856.928400000001
$ws.Range("N122").Value = -16958.4208

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1003.01886
$ws.Range("I122").Value = 944.06976
$ws.Range("J122").Value = 1256.5
$ws.Range("K122").Value = 2832.20928
$ws.Range("L122").Value = 3769.5
$ws.Range("M122").Value = -382.20928
$ws.Range("N122").Value = -8669.5
$ws.Range("H132").Value = 2979.8572
$ws.Range("I132").Value = 2625.244
$ws.Range("K132").Value = 7875.732
$ws.Range("M132").Value = -5345.732
